# "Util updated, duplicate names handled"
#
# The duplicate-name handling fix made every previously-open issue on the
# "Issues" sheet pass testing: each issue's Status becomes "Fixed" and its
# Comments becomes "Tested" (issue #7, "'Cumulative score' is hanging",
# is untouched/unresolved and keeps blank Status/Comments).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

$fixedRows = @(2,3,4,5,6,7,9,10,11,12,13)
foreach ($r in $fixedRows) {
    $ws.Cells.Item($r, 5).Value = "Fixed"    # column E - Status
    $ws.Cells.Item($r, 6).Value = "Tested"   # column F - Comments
}

# Move the saved selection on the Issues sheet from B16 to E16.
$ws.Range("E16").Select()

# Record the (saved) window x-position shift alongside the selection change.
$excel.ActiveWindow.Left = 88.8
